# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 02:22"

# 2) Swap the Albania / Niger rows (shared-string reorder in the source
#    diff manifests as the two country labels trading places, each one
#    keeping its own up-to-date figures)
$ws.Range("A98").Value = "Niger"
$ws.Range("A99").Value = "Albania"

# 3) Refresh the numeric stats that changed with this data pull

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 879430
$ws.Cells.Item(4, 3).Value = 30713
$ws.Cells.Item(4, 5).Value = 744037
$ws.Cells.Item(4, 7).Value = 2110
$ws.Cells.Item(4, 8).Value = 49769

# Row 57 - Argentina
$ws.Cells.Item(57, 2).Value = 3435
$ws.Cells.Item(57, 3).Value = 147
$ws.Cells.Item(57, 5).Value = 2351
$ws.Cells.Item(57, 7).Value = 6
$ws.Cells.Item(57, 8).Value = 165

# Row 98 - now Niger
$ws.Cells.Item(98, 2).Value = 671
$ws.Cells.Item(98, 3).Value = 9
$ws.Cells.Item(98, 4).Value = 256
$ws.Cells.Item(98, 5).Value = 391
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 24

# Row 99 - now Albania
$ws.Cells.Item(99, 2).Value = 663
$ws.Cells.Item(99, 3).Value = 29
$ws.Cells.Item(99, 4).Value = 385
$ws.Cells.Item(99, 5).Value = 251
$ws.Cells.Item(99, 6).Value = 4
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 27

# Row 101 - Burkina Faso
$ws.Cells.Item(101, 2).Value = 616
$ws.Cells.Item(101, 3).Value = 7
$ws.Cells.Item(101, 4).Value = 410
$ws.Cells.Item(101, 5).Value = 165
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(101, 8).Value = 41
